$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.188.92'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.313.65'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.52'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.57'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.99%  '
$ws.Range("E6").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.51%  '
$ws.Range("E7").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E8").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.62%  '
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.51'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.04%  '
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0795'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("E11").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.97%  '
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.03'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +14.32%  '
$ws.Range("E13").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.673.70'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.338.22'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.56%  '
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.819'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.39%  '
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.126.02'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.57'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +10.23%  '
$ws.Range("E19").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("E20").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.16%  '
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.84'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.48'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("E23").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +12.43%  '
$ws.Range("E24").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.83'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.15%  '
$ws.Range("E27").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.93%  '
$ws.Range("E28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.85'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("E29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.12'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("E30").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.04'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("E33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.67'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.23%  '
$ws.Range("E34").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.94%  '
$ws.Range("E35").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.89%  '
$ws.Range("E36").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("E37").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.02%  '
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.80'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.62%  '
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.82'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.61%  '
$ws.Range("E40").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.36'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.000.46'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.66%  '
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.12'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.77%  '
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.74'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.92%  '
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.98'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.69%  '
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.541.26'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("E49").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.31%  '
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.58'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.44%  '
$ws.Range("E51").ClearFormats()
